$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the value in E2 ("!") down to E3, clearing the old E2 cell.
# Use Value2 (rather than the bare Value property) to read the real cell
# content instead of a COM property descriptor.
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("E2").ClearContents()

# Update the active selection to match the authored state.
$ws.Range("F12").Select()
